{"js": "// Locate the paragraph that contains the (originally split) sentence\n// \"Once the user has collected ... Red Pitaya Data Collection\" so the\n// edit does not depend on hard-coded paragraph indices.\nconst results = context.document.body.search(\n  \"Once the user has collected all the necessary data\",\n  { matchCase: false }\n);\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find target paragraph text\");\n}\n\nconst textParagraph = results.items[0].paragraphs.getFirst();\nconst imageParagraph = textParagraph.getNext();\n\n// The paragraph right after the text paragraph holds the inline picture\n// (\"Picture 18\" / ScreenShot5.GIF) that must be removed.\nconst pictures = imageParagraph.inlinePictures;\npictures.load(\"items\");\nawait context.sync();\n\nif (pictures.items.length > 0) {\n  pictures.items[0].delete();\n  await context.sync();\n}\n\n// Merge the two runs of the text paragraph (originally split by the\n// \"_GoBack\" bookmark) into a single run by replacing the paragraph's\n// range with its own (already concatenated) text.\ntextParagraph.load(\"text\");\nawait context.sync();\n\nconst fullText = textParagraph.text;\nconst textRange = textParagraph.getRange();\ntextRange.insertText(fullText, Word.InsertLocation.replace);\nawait context.sync();\n\n// Re-create the \"_GoBack\" bookmark at the start of the now-empty\n// paragraph that used to hold the picture.\nconst imageRange = imageParagraph.getRange();\nimageRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$wdParagraph = 4\n\n$d = $word.ActiveDocument\n\n# --- Step 1: locate the paragraph containing the (originally split) sentence\n# \"Once the user has collected ... Red Pitaya Data Collection\" ---\n$range = $d.Content\n$find = $range.Find\n$find.ClearFormatting()\n$find.Text = \"Once the user has collected all the necessary data\"\n$found = $find.Execute()\n[void]$range.Expand($wdParagraph)\n$textEnd = $range.End\n\n# --- Step 2: find the inline picture (\"Picture 18\" / ScreenShot5.GIF) that sits\n# in the paragraph right after the text paragraph, and delete it ---\n$targetShape = $null\nfor ($i = 1; $i -le $d.InlineShapes.Count; $i++) {\n    $shp = $d.InlineShapes.Item($i)\n    if ($shp.Range.Start -eq $textEnd) {\n        $targetShape = $shp\n    }\n}\nif ($targetShape -ne $null) {\n    $targetShape.Delete()\n}\n\n# --- Step 3: merge the two runs of the text paragraph (originally split by the\n# \"_GoBack\" bookmark) into a single run by replacing the paragraph range's\n# text with its own (already concatenated) text ---\n$range2 = $d.Content\n$find2 = $range2.Find\n$find2.ClearFormatting()\n$find2.Text = \"Once the user has collected all the necessary data\"\n$found2 = $find2.Execute()\n[void]$range2.Expand($wdParagraph)\n$fullText = $range2.Text\n$trimmed = $fullText.Substring(0, $fullText.Length - 1)\n$range2.Text = $trimmed\n\n# --- Step 4: re-create the \"_GoBack\" bookmark at the start of the now-empty\n# paragraph that used to hold the picture ---\n$paras = $d.Paragraphs\n$count = $paras.Count\n$textParaIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $pr = $paras.Item($i).Range\n    if ($pr.Start -le $range2.End -and $pr.End -ge $range2.End) {\n        $textParaIndex = $i\n    }\n}\n$imageParaIndex = $textParaIndex + 1\n$imageParaRange = $paras.Item($imageParaIndex).Range\n$d.Bookmarks.Add(\"_GoBack\", $imageParaRange)\n"}
